# Adds new "Ideas" tasks: a couple of accessibility/input notes, a community
# news-menu idea, a "generalize code for DLC" note, and switches the active
# tab back to the Ideas sheet (selection parked on the newly added rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ideas")

# --- Insert two new rows right after the "Accessibility" bullet list (row 72) ---
$ws.Rows.Item(73).Insert()
$ws.Rows.Item(73).Insert()

$ws.Range("A73").Font.Bold = $true
$ws.Range("A73").HorizontalAlignment = -4108
$ws.Range("B73").Value = "Add a ""Visit forums"" in the help menu"

$ws.Range("A74").Font.Bold = $true
$ws.Range("A74").HorizontalAlignment = -4108
$ws.Range("B74").Value = "Add a fullscreen / windows screen in options"

# --- Insert two new rows inside the "Community" section (after the existing
# bullet list, which now ends at row 92) ---
$ws.Rows.Item(93).Insert()
$ws.Rows.Item(93).Insert()

$ws.Range("A93").Font.Bold = $true
$ws.Range("A93").HorizontalAlignment = -4108
$ws.Range("B93").Value = "News system => make a menu with updates showing like on the iPhone"

$ws.Range("A94").Font.Bold = $true
$ws.Range("A94").HorizontalAlignment = -4108

# --- New bullet under "Core - Other" (row shifted down to 110 by the inserts
# above) ---
$ws.Range("B110").Value = "Generalize my code to support DLC"

# --- Make "Ideas" the active tab again, with the selection parked on the
# first newly-added row ---
$ws.Activate()
$ws.Range("B75").Select()
